# Update countries & provincias Spain
# Applies the daily data refresh: swaps the "Republica Dominicana"/"Sudafrica"
# and "Belice"/"Nueva Caledonia" label order (their ranking rows swapped as
# the underlying dataset moved), refreshes the numeric columns for several
# country rows, and bumps the "datos actualizados" timestamp by one hour.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (A1) ---------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 9 de Mayo de 2020 a las 17:04"

# --- Country label reorder (rows whose ranking position swapped) -------
$ws.Range("A46").Value = "Sudafrica"
$ws.Range("A47").Value = "Republica Dominicana"

$ws.Range("A192").Value = "Nueva Caledonia"
$ws.Range("A193").Value = "Belice"

# --- Row 4: Estados Unidos ----------------------------------------------
$ws.Range("B4").Value = 1325521
$ws.Range("C4").Value = 3736
$ws.Range("D4").Value = 223937
$ws.Range("E4").Value = 1022822
$ws.Range("F4").Value = 16843
$ws.Range("G4").Value = 147
$ws.Range("H4").Value = 78762

# --- Row 19 ---------------------------------------------------------------
$ws.Range("F19").Value = 541

# --- Row 22 ---------------------------------------------------------------
$ws.Range("E22").Value = 2321
$ws.Range("G22").Value = 7
$ws.Range("H22").Value = 1830

# --- Row 46 (now Sudafrica) ------------------------------------------------
$ws.Range("B46").Value = 9420
$ws.Range("C46").Value = 525
$ws.Range("D46").Value = 3983
$ws.Range("E46").Value = 5251
$ws.Range("F46").Value = 77
$ws.Range("G46").Value = 8
$ws.Range("H46").Value = 186

# --- Row 47 (now Republica Dominicana) -------------------------------------
$ws.Range("B47").Value = 9376
$ws.Range("C47").Value = 0
$ws.Range("D47").Value = 2286
$ws.Range("E47").Value = 6710
$ws.Range("F47").Value = 134
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 380

# --- Row 72 -----------------------------------------------------------------
$ws.Range("F72").Value = 29

# --- Row 74 -------------------------------------------------------------------
$ws.Range("B74").Value = 2274
$ws.Range("C74").Value = 7
$ws.Range("D74").Value = 1232
$ws.Range("E74").Value = 934

# --- Row 79 -------------------------------------------------------------------
$ws.Range("B79").Value = 1921
$ws.Range("C79").Value = 49
$ws.Range("E79").Value = 1409
$ws.Range("F79").Value = 56
$ws.Range("G79").Value = 4
$ws.Range("H79").Value = 90

# --- Row 100 ------------------------------------------------------------------
$ws.Range("B100").Value = 892
$ws.Range("C100").Value = 1
$ws.Range("E100").Value = 477

# --- Row 103 ------------------------------------------------------------------
$ws.Range("B103").Value = 847
$ws.Range("C103").Value = 12
$ws.Range("E103").Value = 578

# --- Row 113 ------------------------------------------------------------------
$ws.Range("B113").Value = 692
$ws.Range("C113").Value = 24
$ws.Range("D113").Value = 298
$ws.Range("E113").Value = 357
$ws.Range("G113").Value = 2
$ws.Range("H113").Value = 37

# --- Row 158 ------------------------------------------------------------------
$ws.Range("D158").Value = 64
$ws.Range("E158").Value = 47

# --- Row 192 (now Nueva Caledonia) ---------------------------------------------
$ws.Range("D192").Value = 18
$ws.Range("H192").Value = 0

# --- Row 193 (now Belice) -------------------------------------------------------
$ws.Range("D193").Value = 16
$ws.Range("H193").Value = 2
